{"js": "// Insert \" with:\" immediately after \"Offering 25 years experience\" in the\n// resume summary paragraph, matching the surrounding run formatting\n// (Didactic / Century Gothic, color 343434, size 22 half-points / 11pt).\n\nconst results = context.document.body.search(\"Offering 25 years experience\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Offering 25 years experience\" in document body.');\n}\n\nconst hit = results.items[0];\nconst inserted = hit.insertText(\" with:\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Offering 25 years experience\"\n$find.Execute() | Out-Null\n\n$target = $find.Parent\n$target.Collapse(0)  # wdCollapseEnd\n$target.InsertAfter(\" with:\")\n"}
